$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '9.2'
$ws.Range('B4').NumberFormat = "@"
$ws.Range('B4').Value = 'The Dark Knight'
$ws.Range('C4').NumberFormat = "@"
$ws.Range('C4').Value = '(2008)'
$ws.Range('B5').NumberFormat = "@"
$ws.Range('B5').Value = 'The Godfather: Part II'
$ws.Range('C5').NumberFormat = "@"
$ws.Range('C5').Value = '(1974)'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '9.0'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '8.9'
$ws.Range('B10').NumberFormat = "@"
$ws.Range('B10').Value = 'The Lord of the Rings: The Fellowship of the Ring'
$ws.Range('C10').NumberFormat = "@"
$ws.Range('C10').Value = '(2001)'
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = 'The Good, the Bad and the Ugly'
$ws.Range('C11').NumberFormat = "@"
$ws.Range('C11').Value = '(1966)'
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'Forrest Gump'
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = '(1994)'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.8'
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'Fight Club'
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = '(1999)'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.8'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.7'
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'Se7en'
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = '(1995)'
$ws.Range('B21').NumberFormat = "@"
$ws.Range('B21').Value = 'Seven Samurai'
$ws.Range('C21').NumberFormat = "@"
$ws.Range('C21').Value = '(1954)'
$ws.Range('B22').NumberFormat = "@"
$ws.Range('B22').Value = 'It''s a Wonderful Life'
$ws.Range('C22').NumberFormat = "@"
$ws.Range('C22').Value = '(1946)'
$ws.Range('B23').NumberFormat = "@"
$ws.Range('B23').Value = 'The Silence of the Lambs'
$ws.Range('C23').NumberFormat = "@"
$ws.Range('C23').Value = '(1991)'
$ws.Range('B24').NumberFormat = "@"
$ws.Range('B24').Value = 'Saving Private Ryan'
$ws.Range('C24').NumberFormat = "@"
$ws.Range('C24').Value = '(1998)'
$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'City of God'
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = '(2002)'
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Life Is Beautiful'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = '(1997)'
$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'The Green Mile'
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = '(1999)'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.6'
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'Star Wars'
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = '(1977)'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.6'
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'Interstellar'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = '(2014)'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.6'
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Terminator 2: Judgment Day'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = '(1991)'
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'Back to the Future'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = '(1985)'
$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'Spirited Away'
$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = '(2001)'
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'Psycho'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = '(1960)'
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'Léon: The Professional'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = '(1994)'
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'Parasite'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = '(2019)'
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'The Lion King'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = '(1994)'
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Gladiator'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = '(2000)'
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'American History X'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = '(1998)'
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Spider-Man: No Way Home'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = '(2021)'
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'The Departed'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = '(2006)'
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'The Prestige'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = '(2006)'
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'Casablanca'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = '(1942)'
$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Whiplash'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = '(2014)'
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'The Intouchables'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = '(2011)'
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Modern Times'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = '(1936)'
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Once Upon a Time in the West'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = '(1968)'
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Hara-Kiri'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = '(1962)'
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Grave of the Fireflies'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = '(1988)'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.5'
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Alien'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = '(1979)'
$ws.Range('B53').NumberFormat = "@"
$ws.Range('B53').Value = 'City Lights'
$ws.Range('C53').NumberFormat = "@"
$ws.Range('C53').Value = '(1931)'
$ws.Range('B54').NumberFormat = "@"
$ws.Range('B54').Value = 'Memento'
$ws.Range('C54').NumberFormat = "@"
$ws.Range('C54').Value = '(2000)'
$ws.Range('B56').NumberFormat = "@"
$ws.Range('B56').Value = 'Cinema Paradiso'
$ws.Range('C56').NumberFormat = "@"
$ws.Range('C56').Value = '(1988)'
$ws.Range('B58').NumberFormat = "@"
$ws.Range('B58').Value = 'Django Unchained'
$ws.Range('C58').NumberFormat = "@"
$ws.Range('C58').Value = '(2012)'
$ws.Range('B59').NumberFormat = "@"
$ws.Range('B59').Value = 'WALL·E'
$ws.Range('C59').NumberFormat = "@"
$ws.Range('C59').Value = '(2008)'
$ws.Range('B61').NumberFormat = "@"
$ws.Range('B61').Value = 'Sunset Blvd.'
$ws.Range('C61').NumberFormat = "@"
$ws.Range('C61').Value = '(1950)'
$ws.Range('B62').NumberFormat = "@"
$ws.Range('B62').Value = 'The Shining'
$ws.Range('C62').NumberFormat = "@"
$ws.Range('C62').Value = '(1980)'
$ws.Range('B63').NumberFormat = "@"
$ws.Range('B63').Value = 'Paths of Glory'
$ws.Range('C63').NumberFormat = "@"
$ws.Range('C63').Value = '(1957)'
$ws.Range('B64').NumberFormat = "@"
$ws.Range('B64').Value = 'The Great Dictator'
$ws.Range('C64').NumberFormat = "@"
$ws.Range('C64').Value = '(1940)'
$ws.Range('B65').NumberFormat = "@"
$ws.Range('B65').Value = 'Avengers: Infinity War'
$ws.Range('C65').NumberFormat = "@"
$ws.Range('C65').Value = '(2018)'
$ws.Range('B67').NumberFormat = "@"
$ws.Range('B67').Value = 'The Batman'
$ws.Range('C67').NumberFormat = "@"
$ws.Range('C67').Value = '(2022)'
$ws.Range('B68').NumberFormat = "@"
$ws.Range('B68').Value = 'Aliens'
$ws.Range('C68').NumberFormat = "@"
$ws.Range('C68').Value = '(1986)'
$ws.Range('B69').NumberFormat = "@"
$ws.Range('B69').Value = 'American Beauty'
$ws.Range('C69').NumberFormat = "@"
$ws.Range('C69').Value = '(1999)'
$ws.Range('B70').NumberFormat = "@"
$ws.Range('B70').Value = 'The Dark Knight Rises'
$ws.Range('C70').NumberFormat = "@"
$ws.Range('C70').Value = '(2012)'
$ws.Range('B71').NumberFormat = "@"
$ws.Range('B71').Value = 'Dr. Strangelove or: How I Learned to Stop Worrying and Love the Bomb'
$ws.Range('C71').NumberFormat = "@"
$ws.Range('C71').Value = '(1964)'
$ws.Range('B72').NumberFormat = "@"
$ws.Range('B72').Value = 'Joker'
$ws.Range('C72').NumberFormat = "@"
$ws.Range('C72').Value = '(2019)'
$ws.Range('B73').NumberFormat = "@"
$ws.Range('B73').Value = 'Spider-Man: Into the Spider-Verse'
$ws.Range('C73').NumberFormat = "@"
$ws.Range('C73').Value = '(2018)'
$ws.Range('B74').NumberFormat = "@"
$ws.Range('B74').Value = 'Old Boy'
$ws.Range('C74').NumberFormat = "@"
$ws.Range('C74').Value = '(2003)'
$ws.Range('B75').NumberFormat = "@"
$ws.Range('B75').Value = 'Braveheart'
$ws.Range('C75').NumberFormat = "@"
$ws.Range('C75').Value = '(1995)'
$ws.Range('B76').NumberFormat = "@"
$ws.Range('B76').Value = 'Toy Story'
$ws.Range('C76').NumberFormat = "@"
$ws.Range('C76').Value = '(1995)'
$ws.Range('B77').NumberFormat = "@"
$ws.Range('B77').Value = 'Amadeus'
$ws.Range('C77').NumberFormat = "@"
$ws.Range('C77').Value = '(1984)'
$ws.Range('B78').NumberFormat = "@"
$ws.Range('B78').Value = 'Coco'
$ws.Range('C78').NumberFormat = "@"
$ws.Range('C78').Value = '(2017)'
$ws.Range('B79').NumberFormat = "@"
$ws.Range('B79').Value = 'Inglourious Basterds'
$ws.Range('C79').NumberFormat = "@"
$ws.Range('C79').Value = '(2009)'
$ws.Range('B80').NumberFormat = "@"
$ws.Range('B80').Value = 'The Boat'
$ws.Range('C80').NumberFormat = "@"
$ws.Range('C80').Value = '(1981)'
$ws.Range('B81').NumberFormat = "@"
$ws.Range('B81').Value = 'Avengers: Endgame'
$ws.Range('C81').NumberFormat = "@"
$ws.Range('C81').Value = '(2019)'
$ws.Range('B82').NumberFormat = "@"
$ws.Range('B82').Value = 'Princess Mononoke'
$ws.Range('C82').NumberFormat = "@"
$ws.Range('C82').Value = '(1997)'
$ws.Range('B83').NumberFormat = "@"
$ws.Range('B83').Value = 'Once Upon a Time in America'
$ws.Range('C83').NumberFormat = "@"
$ws.Range('C83').Value = '(1984)'
$ws.Range('B84').NumberFormat = "@"
$ws.Range('B84').Value = 'Good Will Hunting'
$ws.Range('C84').NumberFormat = "@"
$ws.Range('C84').Value = '(1997)'
$ws.Range('B85').NumberFormat = "@"
$ws.Range('B85').Value = 'Toy Story 3'
$ws.Range('C85').NumberFormat = "@"
$ws.Range('C85').Value = '(2010)'
$ws.Range('B86').NumberFormat = "@"
$ws.Range('B86').Value = 'Requiem for a Dream'
$ws.Range('C86').NumberFormat = "@"
$ws.Range('C86').Value = '(2000)'
$ws.Range('B87').NumberFormat = "@"
$ws.Range('B87').Value = '3 Idiots'
$ws.Range('C87').NumberFormat = "@"
$ws.Range('C87').Value = '(2009)'
$ws.Range('B88').NumberFormat = "@"
$ws.Range('B88').Value = 'Singin'' in the Rain'
$ws.Range('C88').NumberFormat = "@"
$ws.Range('C88').Value = '(1952)'
$ws.Range('B89').NumberFormat = "@"
$ws.Range('B89').Value = 'Your Name.'
$ws.Range('C89').NumberFormat = "@"
$ws.Range('C89').Value = '(2016)'
$ws.Range('B91').NumberFormat = "@"
$ws.Range('B91').Value = 'Reservoir Dogs'
$ws.Range('C91').NumberFormat = "@"
$ws.Range('C91').Value = '(1992)'
$ws.Range('B92').NumberFormat = "@"
$ws.Range('B92').Value = 'Eternal Sunshine of the Spotless Mind'
$ws.Range('C92').NumberFormat = "@"
$ws.Range('C92').Value = '(2004)'
$ws.Range('B93').NumberFormat = "@"
$ws.Range('B93').Value = '2001: A Space Odyssey'
$ws.Range('C93').NumberFormat = "@"
$ws.Range('C93').Value = '(1968)'
$ws.Range('B94').NumberFormat = "@"
$ws.Range('B94').Value = 'Citizen Kane'
$ws.Range('C94').NumberFormat = "@"
$ws.Range('C94').Value = '(1941)'
$ws.Range('B95').NumberFormat = "@"
$ws.Range('B95').Value = 'High and Low'
$ws.Range('C95').NumberFormat = "@"
$ws.Range('C95').Value = '(1963)'
$ws.Range('B96').NumberFormat = "@"
$ws.Range('B96').Value = 'Lawrence of Arabia'
$ws.Range('C96').NumberFormat = "@"
$ws.Range('C96').Value = '(1962)'
$ws.Range('B97').NumberFormat = "@"
$ws.Range('B97').Value = 'M'
$ws.Range('C97').NumberFormat = "@"
$ws.Range('C97').Value = '(1931)'
$ws.Range('B98').NumberFormat = "@"
$ws.Range('B98').Value = 'Capernaum'
$ws.Range('C98').NumberFormat = "@"
$ws.Range('C98').Value = '(2018)'
$ws.Range('B99').NumberFormat = "@"
$ws.Range('B99').Value = 'North by Northwest'
$ws.Range('C99').NumberFormat = "@"
$ws.Range('C99').Value = '(1959)'
$ws.Range('B100').NumberFormat = "@"
$ws.Range('B100').Value = 'The Hunt'
$ws.Range('C100').NumberFormat = "@"
$ws.Range('C100').Value = '(2012)'
$ws.Range('B101').NumberFormat = "@"
$ws.Range('B101').Value = 'Vertigo'
$ws.Range('C101').NumberFormat = "@"
$ws.Range('C101').Value = '(1958)'
$ws.Range('B102').NumberFormat = "@"
$ws.Range('B102').Value = 'Amélie'
$ws.Range('C102').NumberFormat = "@"
$ws.Range('C102').Value = '(2001)'
$ws.Range('B103').NumberFormat = "@"
$ws.Range('B103').Value = 'A Clockwork Orange'
$ws.Range('C103').NumberFormat = "@"
$ws.Range('C103').Value = '(1971)'
$ws.Range('B104').NumberFormat = "@"
$ws.Range('B104').Value = 'Full Metal Jacket'
$ws.Range('C104').NumberFormat = "@"
$ws.Range('C104').Value = '(1987)'
$ws.Range('B105').NumberFormat = "@"
$ws.Range('B105').Value = 'Scarface'
$ws.Range('C105').NumberFormat = "@"
$ws.Range('C105').Value = '(1983)'
$ws.Range('B106').NumberFormat = "@"
$ws.Range('B106').Value = 'Double Indemnity'
$ws.Range('C106').NumberFormat = "@"
$ws.Range('C106').Value = '(1944)'
$ws.Range('B107').NumberFormat = "@"
$ws.Range('B107').Value = 'Come and See'
$ws.Range('C107').NumberFormat = "@"
$ws.Range('C107').Value = '(1985)'
$ws.Range('B108').NumberFormat = "@"
$ws.Range('B108').Value = 'The Apartment'
$ws.Range('C108').NumberFormat = "@"
$ws.Range('C108').Value = '(1960)'
$ws.Range('B109').NumberFormat = "@"
$ws.Range('B109').Value = 'Taxi Driver'
$ws.Range('C109').NumberFormat = "@"
$ws.Range('C109').Value = '(1976)'
$ws.Range('B110').NumberFormat = "@"
$ws.Range('B110').Value = 'Hamilton'
$ws.Range('C110').NumberFormat = "@"
$ws.Range('C110').Value = '(2020)'
$ws.Range('B111').NumberFormat = "@"
$ws.Range('B111').Value = 'To Kill a Mockingbird'
$ws.Range('C111').NumberFormat = "@"
$ws.Range('C111').Value = '(1962)'
$ws.Range('B112').NumberFormat = "@"
$ws.Range('B112').Value = 'L.A. Confidential'
$ws.Range('C112').NumberFormat = "@"
$ws.Range('C112').Value = '(1997)'
$ws.Range('B113').NumberFormat = "@"
$ws.Range('B113').Value = 'The Sting'
$ws.Range('C113').NumberFormat = "@"
$ws.Range('C113').Value = '(1973)'
$ws.Range('B114').NumberFormat = "@"
$ws.Range('B114').Value = 'Up'
$ws.Range('C114').NumberFormat = "@"
$ws.Range('C114').Value = '(2009)'
$ws.Range('B115').NumberFormat = "@"
$ws.Range('B115').Value = 'Heat'
$ws.Range('C115').NumberFormat = "@"
$ws.Range('C115').Value = '(1995)'
$ws.Range('B116').NumberFormat = "@"
$ws.Range('B116').Value = 'Snatch'
$ws.Range('C116').NumberFormat = "@"
$ws.Range('C116').Value = '(2000)'
$ws.Range('B117').NumberFormat = "@"
$ws.Range('B117').Value = 'Ikiru'
$ws.Range('C117').NumberFormat = "@"
$ws.Range('C117').Value = '(1952)'
$ws.Range('B118').NumberFormat = "@"
$ws.Range('B118').Value = 'Die Hard'
$ws.Range('C118').NumberFormat = "@"
$ws.Range('C118').Value = '(1988)'
$ws.Range('B119').NumberFormat = "@"
$ws.Range('B119').Value = 'Indiana Jones and the Last Crusade'
$ws.Range('C119').NumberFormat = "@"
$ws.Range('C119').Value = '(1989)'
$ws.Range('B120').NumberFormat = "@"
$ws.Range('B120').Value = 'A Separation'
$ws.Range('C120').NumberFormat = "@"
$ws.Range('C120').Value = '(2011)'
$ws.Range('B121').NumberFormat = "@"
$ws.Range('B121').Value = 'Metropolis'
$ws.Range('C121').NumberFormat = "@"
$ws.Range('C121').Value = '(1927)'
$ws.Range('B122').NumberFormat = "@"
$ws.Range('B122').Value = 'Bicycle Thieves'
$ws.Range('C122').NumberFormat = "@"
$ws.Range('C122').Value = '(1948)'
$ws.Range('B123').NumberFormat = "@"
$ws.Range('B123').Value = 'Incendies'
$ws.Range('C123').NumberFormat = "@"
$ws.Range('C123').Value = '(2010)'
$ws.Range('B124').NumberFormat = "@"
$ws.Range('B124').Value = '1917'
$ws.Range('C124').NumberFormat = "@"
$ws.Range('C124').Value = '(2019)'
$ws.Range('B125').NumberFormat = "@"
$ws.Range('B125').Value = 'Like Stars on Earth'
$ws.Range('C125').NumberFormat = "@"
$ws.Range('C125').Value = '(2007)'
$ws.Range('B126').NumberFormat = "@"
$ws.Range('B126').Value = 'Batman Begins'
$ws.Range('C126').NumberFormat = "@"
$ws.Range('C126').Value = '(2005)'
$ws.Range('B127').NumberFormat = "@"
$ws.Range('B127').Value = 'For a Few Dollars More'
$ws.Range('C127').NumberFormat = "@"
$ws.Range('C127').Value = '(1965)'
$ws.Range('B128').NumberFormat = "@"
$ws.Range('B128').Value = 'Dangal'
$ws.Range('C128').NumberFormat = "@"
$ws.Range('C128').Value = '(2016)'
$ws.Range('B129').NumberFormat = "@"
$ws.Range('B129').Value = 'Downfall'
$ws.Range('C129').NumberFormat = "@"
$ws.Range('C129').Value = '(2004)'
$ws.Range('B130').NumberFormat = "@"
$ws.Range('B130').Value = 'The Kid'
$ws.Range('C130').NumberFormat = "@"
$ws.Range('C130').Value = '(1921)'
$ws.Range('B131').NumberFormat = "@"
$ws.Range('B131').Value = 'Some Like It Hot'
$ws.Range('C131').NumberFormat = "@"
$ws.Range('C131').Value = '(1959)'
$ws.Range('B132').NumberFormat = "@"
$ws.Range('B132').Value = 'The Father'
$ws.Range('C132').NumberFormat = "@"
$ws.Range('C132').Value = '(2020)'
$ws.Range('B133').NumberFormat = "@"
$ws.Range('B133').Value = 'All About Eve'
$ws.Range('C133').NumberFormat = "@"
$ws.Range('C133').Value = '(1950)'
$ws.Range('B134').NumberFormat = "@"
$ws.Range('B134').Value = 'Green Book'
$ws.Range('C134').NumberFormat = "@"
$ws.Range('C134').Value = '(2018)'
$ws.Range('B135').NumberFormat = "@"
$ws.Range('B135').Value = 'The Wolf of Wall Street'
$ws.Range('C135').NumberFormat = "@"
$ws.Range('C135').Value = '(2013)'
$ws.Range('B136').NumberFormat = "@"
$ws.Range('B136').Value = 'Unforgiven'
$ws.Range('C136').NumberFormat = "@"
$ws.Range('C136').Value = '(1992)'
$ws.Range('B137').NumberFormat = "@"
$ws.Range('B137').Value = 'Casino'
$ws.Range('C137').NumberFormat = "@"
$ws.Range('C137').Value = '(1995)'
$ws.Range('B138').NumberFormat = "@"
$ws.Range('B138').Value = 'Pan''s Labyrinth'
$ws.Range('C138').NumberFormat = "@"
$ws.Range('C138').Value = '(2006)'
$ws.Range('B139').NumberFormat = "@"
$ws.Range('B139').Value = 'Judgment at Nuremberg'
$ws.Range('C139').NumberFormat = "@"
$ws.Range('C139').Value = '(1961)'
$ws.Range('B140').NumberFormat = "@"
$ws.Range('B140').Value = 'Ran'
$ws.Range('C140').NumberFormat = "@"
$ws.Range('C140').Value = '(1985)'
$ws.Range('B141').NumberFormat = "@"
$ws.Range('B141').Value = 'A Beautiful Mind'
$ws.Range('C141').NumberFormat = "@"
$ws.Range('C141').Value = '(2001)'
$ws.Range('B142').NumberFormat = "@"
$ws.Range('B142').Value = 'The Sixth Sense'
$ws.Range('C142').NumberFormat = "@"
$ws.Range('C142').Value = '(1999)'
$ws.Range('B143').NumberFormat = "@"
$ws.Range('B143').Value = 'Monty Python and the Holy Grail'
$ws.Range('C143').NumberFormat = "@"
$ws.Range('C143').Value = '(1975)'
$ws.Range('B144').NumberFormat = "@"
$ws.Range('B144').Value = 'There Will Be Blood'
$ws.Range('C144').NumberFormat = "@"
$ws.Range('C144').Value = '(2007)'
$ws.Range('B145').NumberFormat = "@"
$ws.Range('B145').Value = 'The Truman Show'
$ws.Range('C145').NumberFormat = "@"
$ws.Range('C145').Value = '(1998)'
$ws.Range('B146').NumberFormat = "@"
$ws.Range('B146').Value = 'Yojimbo'
$ws.Range('C146').NumberFormat = "@"
$ws.Range('C146').Value = '(1961)'
$ws.Range('B147').NumberFormat = "@"
$ws.Range('B147').Value = 'The Treasure of the Sierra Madre'
$ws.Range('C147').NumberFormat = "@"
$ws.Range('C147').Value = '(1948)'
$ws.Range('B148').NumberFormat = "@"
$ws.Range('B148').Value = 'Shutter Island'
$ws.Range('C148').NumberFormat = "@"
$ws.Range('C148').Value = '(2010)'
$ws.Range('B149').NumberFormat = "@"
$ws.Range('B149').Value = 'The Great Escape'
$ws.Range('C149').NumberFormat = "@"
$ws.Range('C149').Value = '(1963)'
$ws.Range('D149').NumberFormat = "@"
$ws.Range('D149').Value = '8.2'
$ws.Range('B150').NumberFormat = "@"
$ws.Range('B150').Value = 'Rashomon'
$ws.Range('C150').NumberFormat = "@"
$ws.Range('C150').Value = '(1950)'
$ws.Range('B151').NumberFormat = "@"
$ws.Range('B151').Value = 'Jurassic Park'
$ws.Range('C151').NumberFormat = "@"
$ws.Range('C151').Value = '(1993)'
$ws.Range('B152').NumberFormat = "@"
$ws.Range('B152').Value = 'Kill Bill: Vol. 1'
$ws.Range('C152').NumberFormat = "@"
$ws.Range('C152').Value = '(2003)'
$ws.Range('B153').NumberFormat = "@"
$ws.Range('B153').Value = 'Finding Nemo'
$ws.Range('C153').NumberFormat = "@"
$ws.Range('C153').Value = '(2003)'
$ws.Range('B154').NumberFormat = "@"
$ws.Range('B154').Value = 'No Country for Old Men'
$ws.Range('C154').NumberFormat = "@"
$ws.Range('C154').Value = '(2007)'
$ws.Range('B155').NumberFormat = "@"
$ws.Range('B155').Value = 'Raging Bull'
$ws.Range('C155').NumberFormat = "@"
$ws.Range('C155').Value = '(1980)'
$ws.Range('B156').NumberFormat = "@"
$ws.Range('B156').Value = 'The Elephant Man'
$ws.Range('C156').NumberFormat = "@"
$ws.Range('C156').Value = '(1980)'
$ws.Range('B157').NumberFormat = "@"
$ws.Range('B157').Value = 'V for Vendetta'
$ws.Range('C157').NumberFormat = "@"
$ws.Range('C157').Value = '(2005)'
$ws.Range('B158').NumberFormat = "@"
$ws.Range('B158').Value = 'Gone with the Wind'
$ws.Range('C158').NumberFormat = "@"
$ws.Range('C158').Value = '(1939)'
$ws.Range('B159').NumberFormat = "@"
$ws.Range('B159').Value = 'Chinatown'
$ws.Range('C159').NumberFormat = "@"
$ws.Range('C159').Value = '(1974)'
$ws.Range('B160').NumberFormat = "@"
$ws.Range('B160').Value = 'Inside Out'
$ws.Range('C160').NumberFormat = "@"
$ws.Range('C160').Value = '(2015)'
$ws.Range('B161').NumberFormat = "@"
$ws.Range('B161').Value = 'Lock, Stock and Two Smoking Barrels'
$ws.Range('C161').NumberFormat = "@"
$ws.Range('C161').Value = '(1998)'
$ws.Range('B162').NumberFormat = "@"
$ws.Range('B162').Value = 'The Thing'
$ws.Range('C162').NumberFormat = "@"
$ws.Range('C162').Value = '(1982)'
$ws.Range('B163').NumberFormat = "@"
$ws.Range('B163').Value = 'Dial M for Murder'
$ws.Range('C163').NumberFormat = "@"
$ws.Range('C163').Value = '(1954)'
$ws.Range('B164').NumberFormat = "@"
$ws.Range('B164').Value = 'The Secret in Their Eyes'
$ws.Range('C164').NumberFormat = "@"
$ws.Range('C164').Value = '(2009)'
$ws.Range('B165').NumberFormat = "@"
$ws.Range('B165').Value = 'Howl''s Moving Castle'
$ws.Range('C165').NumberFormat = "@"
$ws.Range('C165').Value = '(2004)'
$ws.Range('B166').NumberFormat = "@"
$ws.Range('B166').Value = 'The Bridge on the River Kwai'
$ws.Range('C166').NumberFormat = "@"
$ws.Range('C166').Value = '(1957)'
$ws.Range('B167').NumberFormat = "@"
$ws.Range('B167').Value = 'Trainspotting'
$ws.Range('C167').NumberFormat = "@"
$ws.Range('C167').Value = '(1996)'
$ws.Range('B168').NumberFormat = "@"
$ws.Range('B168').Value = 'Three Billboards Outside Ebbing, Missouri'
$ws.Range('C168').NumberFormat = "@"
$ws.Range('C168').Value = '(2017)'
$ws.Range('B169').NumberFormat = "@"
$ws.Range('B169').Value = 'Warrior'
$ws.Range('C169').NumberFormat = "@"
$ws.Range('C169').Value = '(2011)'
$ws.Range('B170').NumberFormat = "@"
$ws.Range('B170').Value = 'Gran Torino'
$ws.Range('C170').NumberFormat = "@"
$ws.Range('C170').Value = '(2008)'
$ws.Range('B171').NumberFormat = "@"
$ws.Range('B171').Value = 'Fargo'
$ws.Range('C171').NumberFormat = "@"
$ws.Range('C171').Value = '(1996)'
$ws.Range('B172').NumberFormat = "@"
$ws.Range('B172').Value = 'My Neighbor Totoro'
$ws.Range('C172').NumberFormat = "@"
$ws.Range('C172').Value = '(1988)'
$ws.Range('B173').NumberFormat = "@"
$ws.Range('B173').Value = 'Prisoners'
$ws.Range('C173').NumberFormat = "@"
$ws.Range('C173').Value = '(2013)'
$ws.Range('B174').NumberFormat = "@"
$ws.Range('B174').Value = 'Million Dollar Baby'
$ws.Range('C174').NumberFormat = "@"
$ws.Range('C174').Value = '(2004)'
$ws.Range('B175').NumberFormat = "@"
$ws.Range('B175').Value = 'Blade Runner'
$ws.Range('C175').NumberFormat = "@"
$ws.Range('C175').Value = '(1982)'
$ws.Range('B176').NumberFormat = "@"
$ws.Range('B176').Value = 'The Gold Rush'
$ws.Range('C176').NumberFormat = "@"
$ws.Range('C176').Value = '(1925)'
$ws.Range('B177').NumberFormat = "@"
$ws.Range('B177').Value = 'Catch Me If You Can'
$ws.Range('C177').NumberFormat = "@"
$ws.Range('C177').Value = '(2002)'
$ws.Range('B178').NumberFormat = "@"
$ws.Range('B178').Value = 'On the Waterfront'
$ws.Range('C178').NumberFormat = "@"
$ws.Range('C178').Value = '(1954)'
$ws.Range('B179').NumberFormat = "@"
$ws.Range('B179').Value = 'Children of Heaven'
$ws.Range('C179').NumberFormat = "@"
$ws.Range('C179').Value = '(1997)'
$ws.Range('B180').NumberFormat = "@"
$ws.Range('B180').Value = 'The Third Man'
$ws.Range('C180').NumberFormat = "@"
$ws.Range('C180').Value = '(1949)'
$ws.Range('B181').NumberFormat = "@"
$ws.Range('B181').Value = 'Harry Potter and the Deathly Hallows: Part 2'
$ws.Range('C181').NumberFormat = "@"
$ws.Range('C181').Value = '(2011)'
$ws.Range('B182').NumberFormat = "@"
$ws.Range('B182').Value = 'Gone Girl'
$ws.Range('C182').NumberFormat = "@"
$ws.Range('C182').Value = '(2014)'
$ws.Range('B183').NumberFormat = "@"
$ws.Range('B183').Value = 'Ben-Hur'
$ws.Range('C183').NumberFormat = "@"
$ws.Range('C183').Value = '(1959)'
$ws.Range('B184').NumberFormat = "@"
$ws.Range('B184').Value = '12 Years a Slave'
$ws.Range('C184').NumberFormat = "@"
$ws.Range('C184').Value = '(2013)'
$ws.Range('B185').NumberFormat = "@"
$ws.Range('B185').Value = 'The General'
$ws.Range('C185').NumberFormat = "@"
$ws.Range('C185').Value = '(1926)'
$ws.Range('B186').NumberFormat = "@"
$ws.Range('B186').Value = 'The Deer Hunter'
$ws.Range('C186').NumberFormat = "@"
$ws.Range('C186').Value = '(1978)'
$ws.Range('B187').NumberFormat = "@"
$ws.Range('B187').Value = 'Wild Strawberries'
$ws.Range('C187').NumberFormat = "@"
$ws.Range('C187').Value = '(1957)'
$ws.Range('B188').NumberFormat = "@"
$ws.Range('B188').Value = 'Pather Panchali'
$ws.Range('C188').NumberFormat = "@"
$ws.Range('C188').Value = '(1955)'
$ws.Range('B189').NumberFormat = "@"
$ws.Range('B189').Value = 'Before Sunrise'
$ws.Range('C189').NumberFormat = "@"
$ws.Range('C189').Value = '(1995)'
$ws.Range('B190').NumberFormat = "@"
$ws.Range('B190').Value = 'In the Name of the Father'
$ws.Range('C190').NumberFormat = "@"
$ws.Range('C190').Value = '(1993)'
$ws.Range('B191').NumberFormat = "@"
$ws.Range('B191').Value = 'Mr. Smith Goes to Washington'
$ws.Range('C191').NumberFormat = "@"
$ws.Range('C191').Value = '(1939)'
$ws.Range('B192').NumberFormat = "@"
$ws.Range('B192').Value = 'The Grand Budapest Hotel'
$ws.Range('C192').NumberFormat = "@"
$ws.Range('C192').Value = '(2014)'
$ws.Range('B193').NumberFormat = "@"
$ws.Range('B193').Value = 'Room'
$ws.Range('C193').NumberFormat = "@"
$ws.Range('C193').Value = '(2015)'
$ws.Range('B194').NumberFormat = "@"
$ws.Range('B194').Value = 'Sherlock Jr.'
$ws.Range('C194').NumberFormat = "@"
$ws.Range('C194').Value = '(1924)'
$ws.Range('B195').NumberFormat = "@"
$ws.Range('B195').Value = 'Hacksaw Ridge'
$ws.Range('C195').NumberFormat = "@"
$ws.Range('C195').Value = '(2016)'
$ws.Range('B196').NumberFormat = "@"
$ws.Range('B196').Value = 'How to Train Your Dragon'
$ws.Range('C196').NumberFormat = "@"
$ws.Range('C196').Value = '(2010)'
$ws.Range('B197').NumberFormat = "@"
$ws.Range('B197').Value = 'The Wages of Fear'
$ws.Range('C197').NumberFormat = "@"
$ws.Range('C197').Value = '(1953)'
$ws.Range('B198').NumberFormat = "@"
$ws.Range('B198').Value = 'Memories of Murder'
$ws.Range('C198').NumberFormat = "@"
$ws.Range('C198').Value = '(2003)'
$ws.Range('B199').NumberFormat = "@"
$ws.Range('B199').Value = 'The Seventh Seal'
$ws.Range('C199').NumberFormat = "@"
$ws.Range('C199').Value = '(1957)'
$ws.Range('B200').NumberFormat = "@"
$ws.Range('B200').Value = 'Barry Lyndon'
$ws.Range('C200').NumberFormat = "@"
$ws.Range('C200').Value = '(1975)'
$ws.Range('B201').NumberFormat = "@"
$ws.Range('B201').Value = 'The Big Lebowski'
$ws.Range('C201').NumberFormat = "@"
$ws.Range('C201').Value = '(1998)'
$ws.Range('B202').NumberFormat = "@"
$ws.Range('B202').Value = 'Mad Max: Fury Road'
$ws.Range('B203').NumberFormat = "@"
$ws.Range('B203').Value = 'Klaus'
$ws.Range('C203').NumberFormat = "@"
$ws.Range('C203').Value = '(2019)'
$ws.Range('B204').NumberFormat = "@"
$ws.Range('B204').Value = 'Wild Tales'
$ws.Range('C204').NumberFormat = "@"
$ws.Range('C204').Value = '(2014)'
$ws.Range('B205').NumberFormat = "@"
$ws.Range('B205').Value = 'Monsters, Inc.'
$ws.Range('C205').NumberFormat = "@"
$ws.Range('C205').Value = '(2001)'
$ws.Range('B206').NumberFormat = "@"
$ws.Range('B206').Value = 'Mary and Max'
$ws.Range('C206').NumberFormat = "@"
$ws.Range('C206').Value = '(2009)'
$ws.Range('B207').NumberFormat = "@"
$ws.Range('B207').Value = 'Jaws'
$ws.Range('C207').NumberFormat = "@"
$ws.Range('C207').Value = '(1975)'
$ws.Range('B208').NumberFormat = "@"
$ws.Range('B208').Value = 'The Passion of Joan of Arc'
$ws.Range('C208').NumberFormat = "@"
$ws.Range('C208').Value = '(1928)'
$ws.Range('B209').NumberFormat = "@"
$ws.Range('B209').Value = 'Hotel Rwanda'
$ws.Range('C209').NumberFormat = "@"
$ws.Range('C209').Value = '(2004)'
$ws.Range('B210').NumberFormat = "@"
$ws.Range('B210').Value = 'Rocky'
$ws.Range('C210').NumberFormat = "@"
$ws.Range('C210').Value = '(1976)'
$ws.Range('B212').NumberFormat = "@"
$ws.Range('B212').Value = 'Tokyo Story'
$ws.Range('C212').NumberFormat = "@"
$ws.Range('C212').Value = '(1953)'
$ws.Range('B213').NumberFormat = "@"
$ws.Range('B213').Value = 'Platoon'
$ws.Range('C213').NumberFormat = "@"
$ws.Range('C213').Value = '(1986)'
$ws.Range('B214').NumberFormat = "@"
$ws.Range('B214').Value = 'The Terminator'
$ws.Range('C214').NumberFormat = "@"
$ws.Range('C214').Value = '(1984)'
$ws.Range('B215').NumberFormat = "@"
$ws.Range('B215').Value = 'Ford v Ferrari'
$ws.Range('C215').NumberFormat = "@"
$ws.Range('C215').Value = '(2019)'
$ws.Range('B216').NumberFormat = "@"
$ws.Range('B216').Value = 'Stand by Me'
$ws.Range('C216').NumberFormat = "@"
$ws.Range('C216').Value = '(1986)'
$ws.Range('B217').NumberFormat = "@"
$ws.Range('B217').Value = 'Rush'
$ws.Range('C217').NumberFormat = "@"
$ws.Range('C217').Value = '(2013)'
$ws.Range('D217').NumberFormat = "@"
$ws.Range('D217').Value = '8.0'
$ws.Range('B218').NumberFormat = "@"
$ws.Range('B218').Value = 'Into the Wild'
$ws.Range('C218').NumberFormat = "@"
$ws.Range('C218').Value = '(2007)'
$ws.Range('D218').NumberFormat = "@"
$ws.Range('D218').Value = '8.0'
$ws.Range('B219').NumberFormat = "@"
$ws.Range('B219').Value = 'The Wizard of Oz'
$ws.Range('C219').NumberFormat = "@"
$ws.Range('C219').Value = '(1939)'
$ws.Range('D219').NumberFormat = "@"
$ws.Range('D219').Value = '8.0'
$ws.Range('D220').NumberFormat = "@"
$ws.Range('D220').Value = '8.0'
$ws.Range('B221').NumberFormat = "@"
$ws.Range('B221').Value = 'Spotlight'
$ws.Range('C221').NumberFormat = "@"
$ws.Range('C221').Value = '(2015)'
$ws.Range('D221').NumberFormat = "@"
$ws.Range('D221').Value = '8.0'
$ws.Range('B222').NumberFormat = "@"
$ws.Range('B222').Value = 'Network'
$ws.Range('C222').NumberFormat = "@"
$ws.Range('C222').Value = '(1976)'
$ws.Range('D222').NumberFormat = "@"
$ws.Range('D222').Value = '8.0'
$ws.Range('B223').NumberFormat = "@"
$ws.Range('B223').Value = 'Groundhog Day'
$ws.Range('C223').NumberFormat = "@"
$ws.Range('C223').Value = '(1993)'
$ws.Range('D223').NumberFormat = "@"
$ws.Range('D223').Value = '8.0'
$ws.Range('B224').NumberFormat = "@"
$ws.Range('B224').Value = 'The Exorcist'
$ws.Range('C224').NumberFormat = "@"
$ws.Range('C224').Value = '(1973)'
$ws.Range('D224').NumberFormat = "@"
$ws.Range('D224').Value = '8.0'
$ws.Range('B225').NumberFormat = "@"
$ws.Range('B225').Value = 'Ratatouille'
$ws.Range('C225').NumberFormat = "@"
$ws.Range('C225').Value = '(2007)'
$ws.Range('D225').NumberFormat = "@"
$ws.Range('D225').Value = '8.0'
$ws.Range('B227').NumberFormat = "@"
$ws.Range('B227').Value = 'The Incredibles'
$ws.Range('C227').NumberFormat = "@"
$ws.Range('C227').Value = '(2004)'
$ws.Range('B228').NumberFormat = "@"
$ws.Range('B228').Value = 'Dersu Uzala'
$ws.Range('C228').NumberFormat = "@"
$ws.Range('C228').Value = '(1975)'
$ws.Range('B229').NumberFormat = "@"
$ws.Range('B229').Value = 'The Best Years of Our Lives'
$ws.Range('C229').NumberFormat = "@"
$ws.Range('C229').Value = '(1946)'
$ws.Range('B230').NumberFormat = "@"
$ws.Range('B230').Value = 'Before Sunset'
$ws.Range('C230').NumberFormat = "@"
$ws.Range('C230').Value = '(2004)'
$ws.Range('B231').NumberFormat = "@"
$ws.Range('B231').Value = 'Dune'
$ws.Range('C231').NumberFormat = "@"
$ws.Range('C231').Value = '(2021)'
$ws.Range('B233').NumberFormat = "@"
$ws.Range('B233').Value = 'My Father and My Son'
$ws.Range('C233').NumberFormat = "@"
$ws.Range('C233').Value = '(2005)'
$ws.Range('B234').NumberFormat = "@"
$ws.Range('B234').Value = 'The Grapes of Wrath'
$ws.Range('C234').NumberFormat = "@"
$ws.Range('C234').Value = '(1940)'
$ws.Range('B235').NumberFormat = "@"
$ws.Range('B235').Value = 'Cool Hand Luke'
$ws.Range('C235').NumberFormat = "@"
$ws.Range('C235').Value = '(1967)'
$ws.Range('B236').NumberFormat = "@"
$ws.Range('B236').Value = 'To Be or Not to Be'
$ws.Range('C236').NumberFormat = "@"
$ws.Range('C236').Value = '(1942)'
$ws.Range('B237').NumberFormat = "@"
$ws.Range('B237').Value = 'The Battle of Algiers'
$ws.Range('C237').NumberFormat = "@"
$ws.Range('C237').Value = '(1966)'
$ws.Range('B238').NumberFormat = "@"
$ws.Range('B238').Value = 'Amores perros'
$ws.Range('C238').NumberFormat = "@"
$ws.Range('C238').Value = '(2000)'
$ws.Range('B239').NumberFormat = "@"
$ws.Range('B239').Value = 'Pirates of the Caribbean: The Curse of the Black Pearl'
$ws.Range('C239').NumberFormat = "@"
$ws.Range('C239').Value = '(2003)'
$ws.Range('B240').NumberFormat = "@"
$ws.Range('B240').Value = 'The Sound of Music'
$ws.Range('C240').NumberFormat = "@"
$ws.Range('C240').Value = '(1965)'
$ws.Range('B241').NumberFormat = "@"
$ws.Range('B241').Value = 'Life of Brian'
$ws.Range('C241').NumberFormat = "@"
$ws.Range('C241').Value = '(1979)'
$ws.Range('B242').NumberFormat = "@"
$ws.Range('B242').Value = 'The 400 Blows'
$ws.Range('C242').NumberFormat = "@"
$ws.Range('C242').Value = '(1959)'
$ws.Range('B243').NumberFormat = "@"
$ws.Range('B243').Value = 'Persona'
$ws.Range('C243').NumberFormat = "@"
$ws.Range('C243').Value = '(1966)'
$ws.Range('B244').NumberFormat = "@"
$ws.Range('B244').Value = 'It Happened One Night'
$ws.Range('C244').NumberFormat = "@"
$ws.Range('C244').Value = '(1934)'
$ws.Range('B245').NumberFormat = "@"
$ws.Range('B245').Value = 'La Haine'
$ws.Range('C245').NumberFormat = "@"
$ws.Range('C245').Value = '(1995)'
$ws.Range('B246').NumberFormat = "@"
$ws.Range('B246').Value = 'Aladdin'
$ws.Range('C246').NumberFormat = "@"
$ws.Range('C246').Value = '(1992)'
$ws.Range('B247').NumberFormat = "@"
$ws.Range('B247').Value = 'Jai Bhim'
$ws.Range('C247').NumberFormat = "@"
$ws.Range('C247').Value = '(2021)'
$ws.Range('B248').NumberFormat = "@"
$ws.Range('B248').Value = 'Beauty and the Beast'
$ws.Range('C248').NumberFormat = "@"
$ws.Range('C248').Value = '(1991)'
$ws.Range('B249').NumberFormat = "@"
$ws.Range('B249').Value = 'Gandhi'
$ws.Range('B250').NumberFormat = "@"
$ws.Range('B250').Value = 'The Help'
$ws.Range('C250').NumberFormat = "@"
$ws.Range('C250').Value = '(2011)'
$ws.Range('B251').NumberFormat = "@"
$ws.Range('B251').Value = 'The Handmaiden'
$ws.Range('C251').NumberFormat = "@"
$ws.Range('C251').Value = '(2016)'
